$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-09 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-10 Friday", 2) | Out-Null
$d.Content.Find.Execute("533÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "491÷7=", 2) | Out-Null
$d.Content.Find.Execute("527÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "190÷9=", 2) | Out-Null
$d.Content.Find.Execute("843÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "799÷4=", 2) | Out-Null
$d.Content.Find.Execute("750÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "113÷5=", 2) | Out-Null
$d.Content.Find.Execute("315÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "165÷4=", 2) | Out-Null
$d.Content.Find.Execute("677÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "930÷3=", 2) | Out-Null
$d.Content.Find.Execute("844÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "582÷2=", 2) | Out-Null
$d.Content.Find.Execute("822÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "530÷2=", 2) | Out-Null
$d.Content.Find.Execute("949÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "633÷6=", 2) | Out-Null
$d.Content.Find.Execute("280÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "306÷5=", 2) | Out-Null
$d.Content.Find.Execute("816÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "826÷9=", 2) | Out-Null
$d.Content.Find.Execute("161÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "480÷2=", 2) | Out-Null
$d.Content.Find.Execute("545÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "347÷4=", 2) | Out-Null
$d.Content.Find.Execute("150÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "271÷3=", 2) | Out-Null
$d.Content.Find.Execute("314÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "204÷6=", 2) | Out-Null
$d.Content.Find.Execute("996÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "741÷7=", 2) | Out-Null
$d.Content.Find.Execute("425÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "654÷6=", 2) | Out-Null
$d.Content.Find.Execute("630÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "490÷6=", 2) | Out-Null
$d.Content.Find.Execute("663÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "660÷8=", 2) | Out-Null
$d.Content.Find.Execute("743÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "604÷3=", 2) | Out-Null
$d.Content.Find.Execute("776÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "620÷4=", 2) | Out-Null
$d.Content.Find.Execute("406÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "646÷8=", 2) | Out-Null
$d.Content.Find.Execute("850÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "123÷4=", 2) | Out-Null
$d.Content.Find.Execute("607÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "574÷6=", 2) | Out-Null
$d.Content.Find.Execute("199÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "331÷3=", 2) | Out-Null
